$d = $word.ActiveDocument

# =====================================================================
# Part 1: remove the standalone title paragraph
# ("Introdução à programação e pensamento computacional") and add the
# hidden "_GoBack" bookmark (id 0) right at the start of what becomes
# the new first paragraph ("Tipologias e variáveis").
# =====================================================================

$titlePara = $d.Paragraphs(1)
$titlePara.Range.Delete()

# Adding a Bookmarks.Add() collapsed exactly at absolute position 0 is
# mishandled by this host (the bookmarkEnd gets reseated past the next
# paragraph boundary instead of staying collapsed). Work around it by
# inserting a throwaway placeholder character at position 0, anchoring
# the bookmark at the boundary right after it (which behaves correctly),
# then deleting the placeholder again.
$placeholder = $d.Range(0, 0)
$placeholder.InsertBefore("X")
$goBackRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $goBackRange)
$d.Range(0, 1).Delete()

# =====================================================================
# Part 2: "São blocos de instruções que realizam tarefas específ" +
# bookmark "_GoBack" + "icas"  ->  a single run
# "São blocos de instruções que realizam tarefas específicas", with the
# (second) "_GoBack" bookmark removed. The runs that follow ("
# identificados por nomes e parâmetros (assinatura da função)", ". ",
# etc.) must stay exactly as separate runs, untouched.
# =====================================================================

$searchRange = $d.Content.Duplicate
$searchRange.Find.Execute("específ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$beforeBookmarkEnd = $searchRange.End

# Guard the text that must stay untouched from this host's "merge every
# same-look run in the paragraph together" behaviour on edit, by
# dropping a temporary bookmark right after "...específ" + "icas" (i.e.
# right before " identificados..."), so any merge cascades stop there.
$guardRange = $d.Range($beforeBookmarkEnd + 4, $beforeBookmarkEnd + 4)
$d.Bookmarks.Add("zzGuard", $guardRange)

# Remove the real "_GoBack" bookmark sitting between "específ" and
# "icas" - Bookmark.Delete() is surgical and does not disturb the runs.
$d.Bookmarks("_GoBack").Delete()

# Force the "específ" + "icas" runs to actually merge into one run by
# performing a trivial delete+reinsert of the last character of
# "específ" (this is the smallest possible edit touching that run).
$lastCharOfEspecif = $d.Range($beforeBookmarkEnd - 1, $beforeBookmarkEnd)
$lastCharOfEspecif.Delete()
$reinsertPoint = $d.Range($beforeBookmarkEnd - 1, $beforeBookmarkEnd - 1)
$reinsertPoint.InsertBefore("f")

# Drop the temporary guard bookmark again (surgical, no side effects).
$d.Bookmarks("zzGuard").Delete()

Write-Output "edit applied"
